{"js": "// Update benchmark stats table values per the commit:\n// \"Fixed README.md stats and docx preparation for all Renaissance -\n//  JDK 17 - Shenandoah GC tests\"\n//\n// The document is a single-column table; each row holds one stat value.\n// We update the affected rows (by 0-based row index) to their new values,\n// and collapse the three multi-run \"summary\" rows (which had tab-separated\n// numbers) down to the single headline value.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Map of row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"2362\",\n  4: \"0.00002\",\n  5: \"0.03109\",\n  6: \"0.00531\",\n  7: \"0.00010\",\n  8: \"0.03109\",\n  9: \"0.03109\",\n  10: \"0.03109\",\n  11: \"0.47125\",\n  43: \"99.94\",\n  44: \"0.47\",\n  45: \"785\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const idx = parseInt(rowIndex, 10);\n  const cell = table.getCell(idx, 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table values per the commit:\n# \"Fixed README.md stats and docx preparation for all Renaissance -\n#  JDK 17 - Shenandoah GC tests\"\n#\n# The document is a single-column table; each row holds one stat value.\n# We update the affected rows (1-based row index, matching Table.Cell)\n# to their new values. Rows 44/45/46 previously held tab-separated\n# \"raw per-iteration\" numbers across many runs in a single cell; they\n# collapse down to just the single headline value.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$updates = [ordered]@{\n  1  = \"0M\"\n  2  = \"0M\"\n  3  = \"0M\"\n  4  = \"2362\"\n  5  = \"0.00002\"\n  6  = \"0.03109\"\n  7  = \"0.00531\"\n  8  = \"0.00010\"\n  9  = \"0.03109\"\n  10 = \"0.03109\"\n  11 = \"0.03109\"\n  12 = \"0.47125\"\n  44 = \"99.94\"\n  45 = \"0.47\"\n  46 = \"785\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n  $cell = $tbl.Cell($rowIndex, 1)\n  $cell.Range.Text = $updates[$rowIndex]\n}\n"}
